# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# to reflect refreshed cryptocurrency data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are plain-text price strings (e.g. "26.915.68", "22.22")
# that must stay text rather than being auto-coerced into numbers by Excel.
# Prefixing with a literal apostrophe forces text entry; ClearFormats()
# afterwards drops the transient "quote prefix" cell style it introduces so
# the cell keeps the workbook default style, matching the original format.
function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $value
    $cell.ClearFormats()
}

Set-TextValue "D2" "26.915.68"
$ws.Range("E2").Value = "  -0.14%  "
Set-TextValue "D3" "1.549.61"
$ws.Range("E3").Value = "  -0.44%  "
Set-TextValue "D5" "206.29"
$ws.Range("E5").Value = "  -0.31%  "
Set-TextValue "D6" "0.487"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  -0.37%  "
Set-TextValue "D8" "22.22"
$ws.Range("E8").Value = "  +3.36%  "
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("E10").Value = "  +0.57%  "
Set-TextValue "D11" "0.0854"
$ws.Range("E11").Value = "  -0.55%  "
Set-TextValue "D12" "1.770.40"
$ws.Range("E12").Value = "  -0.40%  "
Set-TextValue "D13" "1.548.74"
$ws.Range("E15").Value = "  +0.60%  "
Set-TextValue "D16" "26.903.66"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("E17").Value = "  -0.34%  "
Set-TextValue "D18" "217.00"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("E19").Value = "  +1.51%  "
Set-TextValue "D20" "7.26"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("E24").Value = "  -0.56%  "
Set-TextValue "D25" "154.10"
$ws.Range("E25").Value = "  +0.37%  "
Set-TextValue "D26" "6.62"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("E29").Value = "  -0.39%  "
Set-TextValue "D30" "0.0467"
$ws.Range("E30").Value = "  +1.63%  "
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("E32").Value = "  -0.48%  "
Set-TextValue "D33" "1.416.46"
$ws.Range("E33").Value = "  +3.36%  "
Set-TextValue "D34" "3.07"
$ws.Range("E34").Value = "  +3.63%  "
$ws.Range("E35").Value = "  +1.80%  "
Set-TextValue "D36" "0.967"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +0.25%  "
Set-TextValue "D39" "0.523"
$ws.Range("E39").Value = "  +0.77%  "
Set-TextValue "D40" "0.808"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("E42").Value = "  +3.55%  "
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("E44").Value = "  +1.77%  "
Set-TextValue "D45" "64.51"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("E46").Value = "  +1.20%  "
Set-TextValue "D47" "1.683.79"
$ws.Range("E47").Value = "  -0.40%  "
Set-TextValue "D48" "87.45"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("E49").Value = "  +3.58%  "
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("E51").Value = "  +0.13%  "
